$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header fields: POSITION (B3) and STATUS (B4) ---
# Shared strings are appended in this order so they land on the same
# indices the source workbook used (465=ADMIN AIDE IV, 466=PERMANENT,
# 467=2024).
$ws.Range("B3").Value = "ADMIN AIDE IV"
$ws.Range("B4").Value = "PERMANENT"

# --- SL earned entries for Nov/Dec 2023 (rows 665 & 668) ---
$ws.Range("C665").Value = 1.25
$ws.Range("C668").Value = 1.25

# --- Row 669 becomes a "2024" year-divider label (text, like the other
#     year-header rows elsewhere in the sheet) instead of a 1/1/2024 date.
$ws.Range("A669").Value = "'2024"
$ws.Range("A669").Font.Bold = $true

# --- Row 670: first entry under the new 2024 divider - a new SL(1-0-0)
#     leave credit dated 1/31/2024, with remarks date 1/3/2024.
$ws.Range("A670").Value = 45322
$ws.Range("B670").Value = "SL(1-0-0)"
$ws.Range("H670").Value = 1
$ws.Range("K670").Value = 45294
# K670 needs the same date display format used by the other REMARKS dates
# in this column (e.g. K665); pull that formatting over explicitly.
$ws.Range("K665").Copy()
$ws.Range("K670").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rows 671-773: every subsequent month-end PERIOD date shifts back
#     one day (first-of-month -> last-of-previous-day pattern change).
for ($r = 671; $r -le 773; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 - 1
}

# --- Append a new trailing row (774) continuing the same monthly
#     pattern, copying row 773's formatting/formula and table membership.
$ws.Range("A773:K773").Copy($ws.Range("A774:K774"))
$ws.Range("A774").Value = 48487
$ws.Range("G774").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A8:K774"))

# --- Reflect the user's final on-screen selection (civil-status dropdown
#     at the top, and the newly extended table area at the bottom).
$ws.Range("F2:G2").Select()
$ws.Range("K671").Select()
